$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NA" region label was ambiguous, so rename it to "North Americas"
# for both of the rows that used it (Canada, USA) so the changed/
# differential rows stand out clearly in the report.
$ws.Range("A7").Value = "North Americas"
$ws.Range("A8").Value = "North Americas"

# Select the cells that were changed, matching the reviewer's focus
# on the differential rows when the workbook is reopened.
$ws.Activate()
$ws.Range("A7:A8").Select()
